$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form_8_Interim")

# Delete entire row 13 - shifts rows 14:186 up to 13:185
$ws.Rows.Item(13).Delete()
